$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("G2").Value = "2016-02-18 02:24:54"
$wsZhCn.Range("G3").Value = "2016-02-18 02:24:54"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("G2").Value = "2016-02-18 02:25:14"
$wsDeDe.Range("G3").Value = "2016-02-18 02:25:14"
